$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1, styled like the other header cells (copy style from D1)
$ws.Range("E1").Value = "Тип обучения"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122) # xlPasteFormats

# Dates must stay plain text ("2022-06-10"), not be auto-converted to date
# serials, so force text format before entry and strip it again afterwards
# to keep the cells unstyled, matching the original workbook's look.
$ws.Range("D2:D3").NumberFormat = "@"

# Update row 2
$ws.Range("C2").Value = "Диспечер"
$ws.Range("D2").Value = "2022-06-10"
$ws.Range("E2").Value = "ПБ1"

# Update row 3
$ws.Range("C3").Value = "Диспечер"
$ws.Range("D3").Value = "2022-06-10"
$ws.Range("E3").Value = "ПБ2"

$ws.Range("D2:D3").ClearFormats()
